$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing price values
$ws.Range("C2").Value = 2500
$ws.Range("D2").Value = 5000
$ws.Range("C3").Value = 2000
$ws.Range("D3").Value = 5000

# Add new row 4
$ws.Range("A4").Value = "keychron k3 pro"
$ws.Range("B4").Value = "v2"
$ws.Range("C4").Value = 500
$ws.Range("D4").Value = 1500

# Update selection to match target
$ws.Range("J8").Select()
